$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 25. This shifts the existing
# rows 25..125 down to 26..126 (carrying their formatting/styles with
# them, notably the date style on column D), and extends the sheet's
# used range to A1:T126 automatically.
$ws.Rows(25).Insert()

# Populate the newly inserted (currently blank) row 25 with the new
# market-day record. Columns A,B,C,E..K hold the same constant
# classification data as every other row in this sheet.
$ws.Range("A25").Value = 5
$ws.Range("B25").Value = "Macroferia Regional de Talca"
$ws.Range("C25").Value = "Maule"
$ws.Range("D25").Value = 44971
$ws.Range("E25").Value = 7
$ws.Range("F25").Value = "Fruta"
$ws.Range("G25").Value = 100101
$ws.Range("H25").Value = "Berries"
$ws.Range("I25").Value = 100101001
$ws.Range("J25").Value = "Arándano (blue)"
$ws.Range("K25").Value = "Sin especificar"
$ws.Range("L25").Value = "Primera"
$ws.Range("M25").Value = 100
$ws.Range("N25").Value = 3000
$ws.Range("O25").Value = 3000
$ws.Range("P25").Value = 3000
$ws.Range("Q25").Value = "$/bandeja 2 kilos"
$ws.Range("R25").Value = "Provincia de Curicó"
$ws.Range("S25").Value = 1500
$ws.Range("T25").Value = 2

# Make sure the date column keeps the same custom date/time number
# format used throughout column D.
$ws.Range("D25").NumberFormat = "YYYY-MM-DD HH:MM:SS"
